$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.542.00'
$ws.Cells.Item(2, 5).Value = '  -1.41%  '
$ws.Cells.Item(3, 4).Value = '2.627.89'
$ws.Cells.Item(3, 5).Value = '  +1.04%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 4).Value = '''534.94'
$ws.Cells.Item(5, 5).Value = '  -0.15%  '
$ws.Cells.Item(6, 4).Value = '''142.81'
$ws.Cells.Item(6, 5).Value = '  +1.10%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '''0.567'
$ws.Cells.Item(8, 5).Value = '  +0.17%  '
$ws.Cells.Item(9, 4).Value = '''7.01'
$ws.Cells.Item(9, 5).Value = '  +8.25%  '
$ws.Cells.Item(10, 5).Value = '  -1.68%  '
$ws.Cells.Item(11, 4).Value = '''0.334'
$ws.Cells.Item(11, 5).Value = '  +0.06%  '
$ws.Cells.Item(12, 4).Value = '''0.134'
$ws.Cells.Item(12, 5).Value = '  +0.83%  '
$ws.Cells.Item(13, 4).Value = '3.098.07'
$ws.Cells.Item(13, 5).Value = '  +1.17%  '
$ws.Cells.Item(14, 4).Value = '58.485.96'
$ws.Cells.Item(14, 5).Value = '  -1.38%  '
$ws.Cells.Item(15, 4).Value = '''20.80'
$ws.Cells.Item(15, 5).Value = '  +0.95%  '
$ws.Cells.Item(16, 4).Value = '2.620.87'
$ws.Cells.Item(16, 5).Value = '  +2.20%  '
$ws.Cells.Item(17, 5).Value = '  -0.98%  '
$ws.Cells.Item(18, 4).Value = '''4.39'
$ws.Cells.Item(18, 5).Value = '  +0.76%  '
$ws.Cells.Item(19, 4).Value = '''334.46'
$ws.Cells.Item(19, 5).Value = '  -1.99%  '
$ws.Cells.Item(20, 4).Value = '''10.15'
$ws.Cells.Item(20, 5).Value = '  +0.58%  '
$ws.Cells.Item(21, 4).Value = '''6.23'
$ws.Cells.Item(21, 5).Value = '  -1.97%  '
$ws.Cells.Item(22, 5).Value = '  -0.08%  '
$ws.Cells.Item(23, 4).Value = '''66.02'
$ws.Cells.Item(23, 5).Value = '  -2.32%  '
$ws.Cells.Item(24, 5).Value = '  +1.71%  '
$ws.Cells.Item(25, 5).Value = '  -0.95%  '
$ws.Cells.Item(26, 5).Value = '  -0.06%  '
$ws.Cells.Item(27, 4).Value = '''7.12'
$ws.Cells.Item(27, 5).Value = '  -1.10%  '
$ws.Cells.Item(28, 4).Value = '0.0₃0736'
$ws.Cells.Item(28, 5).Value = '  -0.87%  '
$ws.Cells.Item(29, 5).Value = '  -0.07%  '
$ws.Cells.Item(30, 4).Value = '''1.64'
$ws.Cells.Item(30, 5).Value = '  -1.21%  '
$ws.Cells.Item(31, 4).Value = '''5.86'
$ws.Cells.Item(31, 5).Value = '  +0.90%  '
$ws.Cells.Item(32, 4).Value = '''18.76'
$ws.Cells.Item(32, 5).Value = '  -0.34%  '
$ws.Cells.Item(33, 4).Value = '''150.63'
$ws.Cells.Item(33, 5).Value = '  +0.49%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = '''1.10'
$ws.Cells.Item(35, 5).Value = '  -0.37%  '
$ws.Cells.Item(36, 2).Value = 'SuiNetwork'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(36, 4).Value = '''0.847'
$ws.Cells.Item(36, 5).Value = '  +2.03%  '
$ws.Cells.Item(37, 2).Value = 'Stacks'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(37, 4).Value = '''1.41'
$ws.Cells.Item(37, 5).Value = '  -3.18%  '
$ws.Cells.Item(38, 2).Value = 'Fetch.AI'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(38, 4).Value = '''0.811'
$ws.Cells.Item(38, 5).Value = '  -1.29%  '
$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).Value = '''3.57'
$ws.Cells.Item(39, 5).Value = '  +1.19%  '
$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).Value = '''280.90'
$ws.Cells.Item(40, 5).Value = '  +3.01%  '
$ws.Cells.Item(41, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(41, 4).Value = '''0.999'
$ws.Cells.Item(41, 5).Value = '  -0.09%  '
$ws.Cells.Item(42, 2).Value = 'Mantle'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(42, 4).Value = '''0.599'
$ws.Cells.Item(42, 5).Value = '  +0.17%  '
$ws.Cells.Item(43, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(43, 4).Value = '''10.70'
$ws.Cells.Item(43, 5).Value = '  -0.53%  '
$ws.Cells.Item(44, 2).Value = 'Hedera'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(44, 4).Value = '''0.0531'
$ws.Cells.Item(44, 5).Value = '  +1.64%  '
$ws.Cells.Item(45, 4).Value = '''19.03'
$ws.Cells.Item(45, 5).Value = '  +3.01%  '
$ws.Cells.Item(46, 2).Value = 'Stellar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(46, 4).Value = '''0.0936'
$ws.Cells.Item(46, 5).Value = '  -1.90%  '
$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).Value = '''0.0224'
$ws.Cells.Item(47, 5).Value = '  +0.96%  '
$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).Value = '1.950.07'
$ws.Cells.Item(48, 5).Value = '  +0.35%  '
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).Value = '''4.46'
$ws.Cells.Item(49, 5).Value = '  -0.94%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(50, 4).Value = '''17.89'
$ws.Cells.Item(50, 5).Value = '  -3.47%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '''112.47'
$ws.Cells.Item(51, 5).Value = '  +1.27%  '
